$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial that was bumped by one day
# (45205 -> 45206) for every data row (rows 2 through 261).
$startRow = 2
$endRow = 261
$col = 3  # Column C

for ($r = $startRow; $r -le $endRow; $r++) {
    $cell = $ws.Cells.Item($r, $col)
    if ($cell.Value2 -eq 45205) {
        $cell.Value2 = 45206
    }
}
